# The sheet previously had an "ID" column (A) that is no longer needed.
# Remove it so that the "Groupe de personnages" / "Nombre de personnages"
# data shifts from columns B/C into columns A/B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()

# Rename the two header cells (now in A1/B1), replacing spaces with
# underscores as in the updated workbook.
$ws.Range("A1").Value = "Groupe_de_personnages"
$ws.Range("B1").Value = "Nombre_de_personnages"

# The header row now wraps onto two lines given the new column width,
# so its height grows from the previous single-line size.
$ws.Rows.Item(1).RowHeight = 28.2

# The total row no longer needs its old extra height now that the
# "ID" column (and its taller header) is gone; let it size naturally.
$ws.Rows.Item(6).AutoFit()

# Reflect the final selection left on the sheet.
[void]$ws.Range("N7").Select()
